$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows (salaru/trading/freelance) down
$ws.Rows.Item(2).Insert()

# Fill new row 2: Salary, 2000, 45893.708333333336
$ws.Cells.Item(2, 1).Value = "Salary"
$ws.Cells.Item(2, 2).Value = 2000
$ws.Cells.Item(2, 3).Value = 45893.708333333336

# Copy the date formatting from row 3 (previously row 2) onto the new C2 cell
$ws.Cells.Item(3, 3).Copy()
$ws.Cells.Item(2, 3).PasteSpecial(-4122)

# Append a new row 6: Salary, 4000, 45761.708333333336
$ws.Cells.Item(6, 1).Value = "Salary"
$ws.Cells.Item(6, 2).Value = 4000
$ws.Cells.Item(6, 3).Value = 45761.708333333336

# Copy the date formatting from row 5 onto the new C6 cell
$ws.Cells.Item(5, 3).Copy()
$ws.Cells.Item(6, 3).PasteSpecial(-4122)

$excel.CutCopyMode = 0
